# Logged Week 17 data and fixed Simulate_Season.py tiebreaking method
$wb = $excel.ActiveWorkbook

# --- OFF sheet: update Home (row 2) target depth totals ---
$wsOff = $wb.Worksheets.Item("OFF")
$wsOff.Range("B2").Value = 274
$wsOff.Range("C2").Value = 203
$wsOff.Range("D2").Value = 69
$wsOff.Range("E2").Value = 22
$wsOff.Range("G2").Value = 3

# --- DEF sheet: update Home (row 2) target depth totals ---
$wsDef = $wb.Worksheets.Item("DEF")
$wsDef.Range("B2").Value = 385
$wsDef.Range("C2").Value = 260
$wsDef.Range("D2").Value = 77
$wsDef.Range("E2").Value = 30

$wb.Save()
